# Update cryptos list with latest prices / 1h volume changes
# (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.453.35"
$ws.Range("E2").Value = "  -2.57%  "
$ws.Range("D3").Value = "2.584.55"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'538.64"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "'143.21"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.581"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").Value = "'6.76"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("E10").Value = "  -3.62%  "
$ws.Range("D11").Value = "'0.138"
$ws.Range("E11").Value = "  +3.43%  "
$ws.Range("D12").Value = "'0.331"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("D13").Value = "3.035.29"
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("D14").Value = "58.366.22"
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").Value = "'20.69"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000134"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.589.73"
$ws.Range("E17").Value = "  -3.35%  "
$ws.Range("D18").Value = "'4.46"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'335.83"
$ws.Range("E19").Value = "  -2.65%  "
$ws.Range("D20").Value = "'10.05"
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("D21").Value = "'6.15"
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'66.92"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").Value = "'0.995"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "'0.157"
$ws.Range("E26").Value = "  -5.67%  "
$ws.Range("D27").Value = "'7.03"
$ws.Range("E27").Value = "  -3.85%  "
$ws.Range("D29").Value = "0.0₃0730"
$ws.Range("E29").Value = "  -3.52%  "
$ws.Range("D30").Value = "'1.64"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("D31").Value = "'155.25"
$ws.Range("E31").Value = "  +3.15%  "
$ws.Range("D32").Value = "'5.92"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("D33").Value = "'18.85"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("D34").Value = "'3.89"
$ws.Range("E34").Value = "  -3.75%  "
$ws.Range("D35").Value = "'36.97"
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("E36").Value = "  -5.02%  "
$ws.Range("D37").Value = "'0.826"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "'0.819"
$ws.Range("E38").Value = "  -3.07%  "
$ws.Range("E39").Value = "  -4.15%  "
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("D41").Value = "'282.63"
$ws.Range("E41").Value = "  -3.27%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "'0.590"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").Value = "'10.63"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("D45").Value = "'0.0534"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").Value = "'0.0941"
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").Value = "'18.48"
$ws.Range("E47").Value = "  -4.81%  "
$ws.Range("D48").Value = "'0.0227"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").Value = "1.913.38"
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("D50").Value = "'17.88"
$ws.Range("E50").Value = "  -3.61%  "
$ws.Range("D51").Value = "'4.41"
$ws.Range("E51").Value = "  -4.20%  "
